# Insert two new rows before current row 766, shifting the existing
# data (old rows 766-875) down to rows 768-877.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("766:767").Insert()

# New row 766 (weekly update - "Primera" quality entry)
$ws.Range("A766").Value = 9
$ws.Range("B766").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C766").Value = "Metropolitana"
$ws.Range("D766").Value = 45077
$ws.Range("E766").Value = 13
$ws.Range("F766").Value = 100114014
$ws.Range("G766").Value = "Betarraga"
$ws.Range("H766").Value = "Sin especificar"
$ws.Range("I766").Value = "Primera"
$ws.Range("J766").Value = 7000
$ws.Range("K766").Value = 90
$ws.Range("L766").Value = 100
$ws.Range("M766").Value = 95
$ws.Range("N766").Value = "`$/unidad"
$ws.Range("O766").Value = "Región Metropolitana"
$ws.Range("P766").Value = 95
$ws.Range("Q766").Value = 1
$ws.Range("R766").Value = "Hortaliza"

# New row 767 (weekly update - "Segunda" quality entry)
$ws.Range("A767").Value = 9
$ws.Range("B767").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C767").Value = "Metropolitana"
$ws.Range("D767").Value = 45077
$ws.Range("E767").Value = 13
$ws.Range("F767").Value = 100114014
$ws.Range("G767").Value = "Betarraga"
$ws.Range("H767").Value = "Sin especificar"
$ws.Range("I767").Value = "Segunda"
$ws.Range("J767").Value = 6100
$ws.Range("K767").Value = 80
$ws.Range("L767").Value = 80
$ws.Range("M767").Value = 80
$ws.Range("N767").Value = "`$/unidad"
$ws.Range("O767").Value = "Región Metropolitana"
$ws.Range("P767").Value = 80
$ws.Range("Q767").Value = 1
$ws.Range("R767").Value = "Hortaliza"
